# Update odds/closing-line data for "Slovenia Prva Liga" sheet and remove the
# trailing (now stale) fixture row, per the 04-04-2024 23:22 base refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 137 (id 135, match 6814753): closing-line refresh ---
$ws.Range("O137").Value = 4.2
$ws.Range("P137").Value = 1.45
$ws.Range("R137").Value = 1.775
$ws.Range("S137").Value = 2.025
$ws.Range("U137").Value = 1.825
$ws.Range("V137").Value = 1.975

# --- Row 138 (id 136, match 8035687): closing-line refresh ---
$ws.Range("U138").Value = 1.85
$ws.Range("V138").Value = 1.95

# --- Row 139 (id 137): now carries the fixture that used to sit on row 140 ---
$ws.Range("B139").Value = 6814435
$ws.Range("E139").Value = 45388.52083333334
$ws.Range("F139").Value = "NK Radomlje"
$ws.Range("G139").Value = "FC Koper"
$ws.Range("K139").Value = 2.55
$ws.Range("L139").Value = 3.25
$ws.Range("M139").Value = 2.55
$ws.Range("N139").Value = 2.45
$ws.Range("O139").Value = 3.25
$ws.Range("P139").Value = 2.625
$ws.Range("Q139").Value = 0
$ws.Range("T139").Value = 2.25
$ws.Range("U139").Value = 1.775
$ws.Range("V139").Value = 2.025

# --- Row 140 (id 138): now carries the fixture that used to sit on row 141 ---
$ws.Range("B140").Value = 6837117
$ws.Range("E140").Value = 45388.63541666666
$ws.Range("F140").Value = "NS Mura"
$ws.Range("G140").Value = "NK Celje"
$ws.Range("K140").Value = 5.25
$ws.Range("L140").Value = 4.2
$ws.Range("M140").Value = 1.5
$ws.Range("N140").Value = 5.25
$ws.Range("O140").Value = 4.2
$ws.Range("P140").Value = 1.5
$ws.Range("Q140").Value = 1
$ws.Range("R140").Value = 1.925
$ws.Range("S140").Value = 1.875
$ws.Range("T140").Value = 2.75
$ws.Range("U140").Value = 1.975
$ws.Range("V140").Value = 1.825

# --- Row 141 (id 139): now carries the fixture that used to sit on row 142 ---
$ws.Range("B141").Value = 6814434
$ws.Range("E141").Value = 45389.41666666666
$ws.Range("F141").Value = "NK Bravo"
$ws.Range("G141").Value = "NK Domzale"
$ws.Range("K141").Value = 1.833
$ws.Range("L141").Value = 3.25
$ws.Range("M141").Value = 4
$ws.Range("N141").Value = 1.833
$ws.Range("O141").Value = 3.25
$ws.Range("P141").Value = 4
$ws.Range("Q141").Value = -0.5
$ws.Range("R141").Value = 1.825
$ws.Range("S141").Value = 1.975
$ws.Range("T141").Value = 2.25
$ws.Range("U141").Value = 1.8
$ws.Range("V141").Value = 2

# --- Row 142 (id 140, match 6814434) is now merged into row 141 above; delete it ---
$ws.Rows(142).Delete() | Out-Null
